# Equal Profit Rate updated.xlsx - apply commit: "Trace for Reset price working
# adequately to debug the process itself"
#
# Changes applied:
#  1. On the "Selected Prices" sheet, cell K22 (Reset price for DI) is changed
#     from 1 to 1.2, and K23 (Reset price for DII) is changed from 1 to 0.9.
#  2. Several dependent formulas are corrected/re-pointed so the "Reset
#     price" trace recalculates correctly:
#       E22: =F$9*C$23   -> =F$9*C$21
#       L22: =E$9*C$16   -> =E$9*C$22
#       R22: =K22/Q$18   -> =K22/Q$24
#       E23: =F$10*C$23  -> =F$10*C$21
#       L23: =E$10*C$16  -> =E$10*C$22
#       R23: =K23/Q$18   -> =K23/Q$24
#  3. The "Output unit prices" note in cell W21 changes from
#     "NOT RIGHT YET!" to "RIGHT" now that the trace checks out.
#  4. A threaded review comment is added to C21 explaining the assumption
#     used for the purpose of the demonstration (price of labour power held
#     constant).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Selected Prices")

# --- 1 & 2: fix up the "Reset price" trace formulas / inputs ------------
$ws.Range("E22").Formula = "=F$9*C$21"
$ws.Range("L22").Formula = "=E$9*C$22"
$ws.Range("R22").Formula = "=K22/Q$24"

$ws.Range("E23").Formula = "=F$10*C$21"
$ws.Range("L23").Formula = "=E$10*C$22"
$ws.Range("R23").Formula = "=K23/Q$24"

$ws.Range("K22").Value = 1.2
$ws.Range("K23").Value = 0.9

# --- 3: flip the status note now that the numbers reconcile -------------
$ws.Range("W21").Value = "RIGHT"

# --- 4: leave a threaded comment documenting the assumption -------------
$commentText = "We assume that the price of labour power remains the same, for the purpose of demonstration. In other words, the real wage rises"

$cell = $ws.Range("C21")
# First add a plain comment - this registers the comment font/box styling
# that Excel uses for notes - then replace it with a threaded comment,
# which is what modern Excel actually stores review comments as.
$legacyComment = $cell.AddComment("placeholder")
$legacyComment.Delete()
$threadedComment = $cell.AddCommentThreaded($commentText)

# Park the selection near the note, mirroring where the author left off.
$ws.Range("W22").Select() | Out-Null
